# Pharma_Society_Report.xlsx update
# 1. Rename the worksheet from "Report" to "Sheet1"
# 2. Correct the membership-count column (B) so the values are stored as
#    real numbers (FLASCO=400, GASCO=500, IOS=50, IOWA=126, MOASC=600)
#    instead of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

$ws.Range("B2").Value = 400
$ws.Range("B3").Value = 500
$ws.Range("B4").Value = 50
$ws.Range("B5").Value = 126
$ws.Range("B6").Value = 600
